# Updated cryptos list on Fri Feb 16 02:59:44 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value to a cell while keeping it text-typed (not
# auto-converted to a number by Excel) and without leaving a residual
# number-format style applied to the cell.
function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# --- Simple value updates (Price / Volume(1h) columns) ---
Set-TextValue "D2" "52.157.25"
$ws.Range("E2").Value = "  +0.12%  "
Set-TextValue "D3" "2.845.02"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("E4").Value = "  -0.04%  "
Set-TextValue "D5" "361.51"
$ws.Range("E5").Value = "  +5.44%  "
Set-TextValue "D6" "114.00"
$ws.Range("E6").Value = "  -2.32%  "
Set-TextValue "D7" "0.571"
$ws.Range("E7").Value = "  +4.93%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +4.21%  "
Set-TextValue "D10" "41.75"
$ws.Range("E10").Value = "  -0.79%  "
Set-TextValue "D11" "0.0864"
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("E12").Value = "  +1.13%  "
Set-TextValue "D13" "20.00"
$ws.Range("E13").Value = "  -0.55%  "
Set-TextValue "D14" "7.78"
$ws.Range("E14").Value = "  +1.92%  "
Set-TextValue "D15" "3.289.42"
$ws.Range("E15").Value = "  +1.79%  "
Set-TextValue "D16" "2.838.06"
$ws.Range("E16").Value = "  +1.31%  "
Set-TextValue "D17" "0.912"
$ws.Range("E17").Value = "  +2.80%  "
Set-TextValue "D18" "52.002.00"
$ws.Range("E18").Value = "  +0.00%  "
Set-TextValue "D19" "7.49"
$ws.Range("E19").Value = "  +7.57%  "
$ws.Range("E20").Value = "  -1.80%  "
Set-TextValue "D21" "13.59"
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("E22").Value = "  +0.41%  "
Set-TextValue "D23" "70.30"
$ws.Range("E23").Value = "  +0.11%  "
Set-TextValue "D24" "267.67"
$ws.Range("E24").Value = "  -3.69%  "
Set-TextValue "D25" "2.83"
$ws.Range("E25").Value = "  -0.56%  "
Set-TextValue "D26" "27.25"
$ws.Range("E26").Value = "  +1.40%  "
Set-TextValue "D27" "0.999"
$ws.Range("E27").Value = "  +0.00%  "
Set-TextValue "D28" "10.45"
$ws.Range("E28").Value = "  +2.57%  "
Set-TextValue "D30" "53.24"
$ws.Range("E30").Value = "  +5.17%  "
$ws.Range("E31").Value = "  -1.71%  "

# --- Row 32 / Row 33 swap (VeChain <-> InjectiveProtocol) ---
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D32" "34.21"
$ws.Range("E32").Value = "  -1.53%  "

$ws.Range("B33").Value = "VeChain"
$ws.Range("C33").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D33" "0.0455"
$ws.Range("E33").Value = "  +23.24%  "

$ws.Range("E34").Value = "  +3.98%  "
$ws.Range("E35").Value = "  +7.98%  "
$ws.Range("E36").Value = "  +2.64%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("E39").Value = "  -2.29%  "
Set-TextValue "D40" "18.39"
$ws.Range("E40").Value = "  -3.05%  "
Set-TextValue "D41" "23.98"
$ws.Range("E41").Value = "  +2.83%  "
$ws.Range("E42").Value = "  +2.19%  "

# --- Row 43 / Row 44 swap (Monero <-> Stacks) ---
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D43" "2.57"
$ws.Range("E43").Value = "  -7.47%  "

$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D44" "128.10"
$ws.Range("E44").Value = "  +0.12%  "

Set-TextValue "D45" "2.27"
$ws.Range("E45").Value = "  -2.98%  "
Set-TextValue "D46" "2.123.91"
$ws.Range("E46").Value = "  +0.84%  "
Set-TextValue "D47" "3.41"
$ws.Range("E47").Value = "  +2.34%  "
$ws.Range("E49").Value = "  +9.51%  "
Set-TextValue "D50" "5.86"
$ws.Range("E50").Value = "  +5.40%  "
$ws.Range("E51").Value = "  +1.16%  "
